$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "b"
$ws.Range("J2").Value = "Acknowledge (Backchannel)"
$ws.Range("I24").Value = "sd"
$ws.Range("J24").Value = "Statement-non-opinion"
$ws.Range("I25").Value = "%"
$ws.Range("J25").Value = "Uninterpretable"
$ws.Range("I28").Value = "sd"
$ws.Range("J28").Value = "Statement-non-opinion"
$ws.Range("I33").Value = "sd"
$ws.Range("J33").Value = "Statement-non-opinion"
$ws.Range("I53").Value = "sd"
$ws.Range("J53").Value = "Statement-non-opinion"
$ws.Range("I54").Value = "sv"
$ws.Range("J54").Value = "Statement-opinion"
$ws.Range("I56").Value = "aa"
$ws.Range("J56").Value = "Agree/Accept"
$ws.Range("I62").Value = "sv"
$ws.Range("J62").Value = "Statement-opinion"
$ws.Range("I63").Value = "sv"
$ws.Range("J63").Value = "Statement-opinion"
$ws.Range("I73").Value = "qy"
$ws.Range("J73").Value = "Yes-No-Question"
$ws.Range("I105").Value = "sv"
$ws.Range("J105").Value = "Statement-opinion"
$ws.Range("I120").Value = "aa"
$ws.Range("J120").Value = "Agree/Accept"
$ws.Range("I138").Value = "aa"
$ws.Range("J138").Value = "Agree/Accept"
$ws.Range("I143").Value = "ba"
$ws.Range("J143").Value = "Appreciation"
$ws.Range("I147").Value = "aa"
$ws.Range("J147").Value = "Agree/Accept"
$ws.Range("I148").Value = "sd"
$ws.Range("J148").Value = "Statement-non-opinion"
$ws.Range("I153").Value = "aa"
$ws.Range("J153").Value = "Agree/Accept"
$ws.Range("I160").Value = "sd"
$ws.Range("J160").Value = "Statement-non-opinion"
$ws.Range("I170").Value = "sd"
$ws.Range("J170").Value = "Statement-non-opinion"
$ws.Range("I172").Value = "sd"
$ws.Range("J172").Value = "Statement-non-opinion"
$ws.Range("I173").Value = "ba"
$ws.Range("J173").Value = "Appreciation"
$ws.Range("I183").Value = "sv"
$ws.Range("J183").Value = "Statement-opinion"
$ws.Range("I188").Value = "qy"
$ws.Range("J188").Value = "Yes-No-Question"
$ws.Range("I189").Value = "aa"
$ws.Range("J189").Value = "Agree/Accept"
$ws.Range("I190").Value = "sv"
$ws.Range("J190").Value = "Statement-opinion"
$ws.Range("I199").Value = "b"
$ws.Range("J199").Value = "Acknowledge (Backchannel)"
$ws.Range("I201").Value = "sd"
$ws.Range("J201").Value = "Statement-non-opinion"
$ws.Range("I205").Value = "%"
$ws.Range("J205").Value = "Uninterpretable"
$ws.Range("I211").Value = "sv"
$ws.Range("J211").Value = "Statement-opinion"
$ws.Range("I218").Value = "sd"
$ws.Range("J218").Value = "Statement-non-opinion"
$ws.Range("I225").Value = "sd"
$ws.Range("J225").Value = "Statement-non-opinion"
$ws.Range("I238").Value = "sv"
$ws.Range("J238").Value = "Statement-opinion"
$ws.Range("I251").Value = "%"
$ws.Range("J251").Value = "Uninterpretable"
$ws.Range("I254").Value = "b"
$ws.Range("J254").Value = "Acknowledge (Backchannel)"
$ws.Range("I274").Value = "qy"
$ws.Range("J274").Value = "Yes-No-Question"
$ws.Range("I288").Value = "b"
$ws.Range("J288").Value = "Acknowledge (Backchannel)"
$ws.Range("I305").Value = "%"
$ws.Range("J305").Value = "Uninterpretable"
$ws.Range("I316").Value = "sd"
$ws.Range("J316").Value = "Statement-non-opinion"
$ws.Range("I317").Value = "ba"
$ws.Range("J317").Value = "Appreciation"
$ws.Range("I325").Value = "aa"
$ws.Range("J325").Value = "Agree/Accept"
$ws.Range("I326").Value = "aa"
$ws.Range("J326").Value = "Agree/Accept"
$ws.Range("I328").Value = "sd"
$ws.Range("J328").Value = "Statement-non-opinion"
$ws.Range("I333").Value = "sd"
$ws.Range("J333").Value = "Statement-non-opinion"
$ws.Range("I334").Value = "sd"
$ws.Range("J334").Value = "Statement-non-opinion"
$ws.Range("I338").Value = "sd"
$ws.Range("J338").Value = "Statement-non-opinion"
$ws.Range("I339").Value = "b"
$ws.Range("J339").Value = "Acknowledge (Backchannel)"
$ws.Range("I345").Value = "ba"
$ws.Range("J345").Value = "Appreciation"
$ws.Range("I346").Value = "aa"
$ws.Range("J346").Value = "Agree/Accept"
$ws.Range("I349").Value = "sd"
$ws.Range("J349").Value = "Statement-non-opinion"
$ws.Range("I356").Value = "aa"
$ws.Range("J356").Value = "Agree/Accept"
$ws.Range("I357").Value = "qy"
$ws.Range("J357").Value = "Yes-No-Question"
$ws.Range("I366").Value = "sv"
$ws.Range("J366").Value = "Statement-opinion"
$ws.Range("I376").Value = "sv"
$ws.Range("J376").Value = "Statement-opinion"
$ws.Range("I379").Value = "sd"
$ws.Range("J379").Value = "Statement-non-opinion"
$ws.Range("I381").Value = "sv"
$ws.Range("J381").Value = "Statement-opinion"
$ws.Range("I406").Value = "aa"
$ws.Range("J406").Value = "Agree/Accept"
$ws.Range("I411").Value = "aa"
$ws.Range("J411").Value = "Agree/Accept"
$ws.Range("I431").Value = "b"
$ws.Range("J431").Value = "Acknowledge (Backchannel)"
$ws.Range("I434").Value = "ba"
$ws.Range("J434").Value = "Appreciation"
$ws.Range("I447").Value = "sd"
$ws.Range("J447").Value = "Statement-non-opinion"
$ws.Range("I458").Value = "sd"
$ws.Range("J458").Value = "Statement-non-opinion"
